$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.141.39'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.057.55'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '250.17'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.670'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.29'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +10.63%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0795'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.18'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +8.74%  '
$ws.Range('D13').Value = '2.356.04'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.827'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.76'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +10.65%  '
$ws.Range('D16').Value = '2.060.07'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.72'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +32.38%  '
$ws.Range('D18').Value = '37.115.93'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '75.60'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.10%  '
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.49'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '239.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +14.00%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.24%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '169.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.31'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.126'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.87'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +6.76%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.14'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +9.79%  '
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.57'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0897'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.99%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  +7.11%  '
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.20'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +27.77%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.12'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +11.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.84'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '98.23'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.09'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.31%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.47'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('D48').Value = '1.295.81'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('D51').Value = '2.239.94'
$ws.Range('E51').Value = '  -0.39%  '
